$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1533
$ws1.Range("F5").Value = 234
$ws1.Range("F7").Value = 161
$ws1.Range("F8").Value = 6302
$ws1.Range("F12").Value = 5324
$ws1.Range("F15").Value = 1196
$ws1.Range("F16").Value = 1
$ws1.Range("F21").Value = 307
$ws1.Range("F23").Value = 3803

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 91

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 91
$ws4.Range("F5").Value = 1533
$ws4.Range("F6").Value = 234
$ws4.Range("F8").Value = 161
$ws4.Range("F9").Value = 6302
$ws4.Range("F13").Value = 5324
$ws4.Range("F16").Value = 1196
$ws4.Range("F17").Value = 1
$ws4.Range("F22").Value = 307
$ws4.Range("F24").Value = 3803
